$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update crypto price (D) and volume-change (E) columns for rows 2-51
# D-column values that look like plain numbers must be forced to text
# (NumberFormat "@" then reset Style to Normal so no stray style survives)

$c = $ws.Range("D2")
$c.NumberFormat = "@"
$c.Value = "67.978.26"
$c.Style = "Normal"
$ws.Range("E2").Value = "  +1.46%  "

$c = $ws.Range("D3")
$c.NumberFormat = "@"
$c.Value = "3.543.09"
$c.Style = "Normal"
$ws.Range("E3").Value = "  +0.35%  "

$c = $ws.Range("D4")
$c.NumberFormat = "@"
$c.Value = "0.999"
$c.Style = "Normal"
$ws.Range("E4").Value = "  -0.20%  "

$c = $ws.Range("D5")
$c.NumberFormat = "@"
$c.Value = "616.16"
$c.Style = "Normal"
$ws.Range("E5").Value = "  +0.74%  "

$c = $ws.Range("D6")
$c.NumberFormat = "@"
$c.Value = "152.83"
$c.Style = "Normal"
$ws.Range("E6").Value = "  -1.19%  "

$c = $ws.Range("D7")
$c.NumberFormat = "@"
$c.Value = "3.543.11"
$c.Style = "Normal"
$ws.Range("E7").Value = "  +0.58%  "

$ws.Range("E8").Value = "  -0.10%  "

$ws.Range("E9").Value = "  -0.43%  "

$c = $ws.Range("D10")
$c.NumberFormat = "@"
$c.Value = "0.141"
$c.Style = "Normal"
$ws.Range("E10").Value = "  -0.90%  "

$ws.Range("E11").Value = "  +3.42%  "

$c = $ws.Range("D12")
$c.NumberFormat = "@"
$c.Value = "0.428"
$c.Style = "Normal"
$ws.Range("E12").Value = "  -1.07%  "

$ws.Range("E13").Value = "  -0.80%  "

$c = $ws.Range("D14")
$c.NumberFormat = "@"
$c.Value = "4.142.10"
$c.Style = "Normal"
$ws.Range("E14").Value = "  +0.29%  "

$c = $ws.Range("D15")
$c.NumberFormat = "@"
$c.Value = "32.16"
$c.Style = "Normal"
$ws.Range("E15").Value = "  +0.54%  "

$c = $ws.Range("D16")
$c.NumberFormat = "@"
$c.Value = "3.533.74"
$c.Style = "Normal"
$ws.Range("E16").Value = "  -0.88%  "

$c = $ws.Range("D17")
$c.NumberFormat = "@"
$c.Value = "67.684.99"
$c.Style = "Normal"
$ws.Range("E17").Value = "  +0.92%  "

$ws.Range("E18").Value = "  -0.65%  "

$c = $ws.Range("D19")
$c.NumberFormat = "@"
$c.Value = "6.42"
$c.Style = "Normal"
$ws.Range("E19").Value = "  +0.44%  "

$c = $ws.Range("D20")
$c.NumberFormat = "@"
$c.Value = "15.40"
$c.Style = "Normal"
$ws.Range("E20").Value = "  -0.28%  "

$ws.Range("E21").Value = "  +3.15%  "

$c = $ws.Range("D22")
$c.NumberFormat = "@"
$c.Value = "448.29"
$c.Style = "Normal"
$ws.Range("E22").Value = "  -0.81%  "

$ws.Range("E23").Value = "  -2.42%  "

$c = $ws.Range("D24")
$c.NumberFormat = "@"
$c.Value = "77.66"
$c.Style = "Normal"
$ws.Range("E24").Value = "  -2.25%  "

$ws.Range("E25").Value = "  +6.24%  "

$c = $ws.Range("D26")
$c.NumberFormat = "@"
$c.Value = "3.683.53"
$c.Style = "Normal"
$ws.Range("E26").Value = "  +0.18%  "

$ws.Range("E27").Value = "  +0.22%  "

$ws.Range("E28").Value = "  -0.65%  "

$c = $ws.Range("D29")
$c.NumberFormat = "@"
$c.Value = "8.70"
$c.Style = "Normal"
$ws.Range("E29").Value = "  +4.16%  "

$ws.Range("E30").Value = "  -0.97%  "

$c = $ws.Range("D31")
$c.NumberFormat = "@"
$c.Value = "1.61"
$c.Style = "Normal"
$ws.Range("E31").Value = "  -4.32%  "

$c = $ws.Range("D32")
$c.NumberFormat = "@"
$c.Value = "0.169"
$c.Style = "Normal"
$ws.Range("E32").Value = "  +7.02%  "

$ws.Range("E33").Value = "  -0.17%  "

$c = $ws.Range("D34")
$c.NumberFormat = "@"
$c.Value = "25.99"
$c.Style = "Normal"
$ws.Range("E34").Value = "  +0.10%  "

$c = $ws.Range("D35")
$c.NumberFormat = "@"
$c.Value = "6.24"
$c.Style = "Normal"
$ws.Range("E35").Value = "  +0.74%  "

$c = $ws.Range("D36")
$c.NumberFormat = "@"
$c.Value = "3.530.14"
$c.Style = "Normal"
$ws.Range("E36").Value = "  -0.02%  "

$ws.Range("E37").Value = "  -1.99%  "

$c = $ws.Range("D38")
$c.NumberFormat = "@"
$c.Value = "8.07"
$c.Style = "Normal"
$ws.Range("E38").Value = "  -0.46%  "

$ws.Range("E39").Value = "  -0.01%  "

$c = $ws.Range("D40")
$c.NumberFormat = "@"
$c.Value = "0.999"
$c.Style = "Normal"
$ws.Range("E40").Value = "  -0.22%  "

$c = $ws.Range("D41")
$c.NumberFormat = "@"
$c.Value = "176.86"
$c.Style = "Normal"
$ws.Range("E41").Value = "  -0.27%  "

$c = $ws.Range("D42")
$c.NumberFormat = "@"
$c.Value = "2.21"
$c.Style = "Normal"
$ws.Range("E42").Value = "  +3.02%  "

$ws.Range("E43").Value = "  +1.95%  "

$c = $ws.Range("D44")
$c.NumberFormat = "@"
$c.Value = "5.44"
$c.Style = "Normal"
$ws.Range("E44").Value = "  -3.15%  "

$c = $ws.Range("D45")
$c.NumberFormat = "@"
$c.Value = "0.888"
$c.Style = "Normal"
$ws.Range("E45").Value = "  -0.48%  "

$c = $ws.Range("D46")
$c.NumberFormat = "@"
$c.Value = "28.96"
$c.Style = "Normal"
$ws.Range("E46").Value = "  +1.66%  "

$c = $ws.Range("D47")
$c.NumberFormat = "@"
$c.Value = "45.51"
$c.Style = "Normal"
$ws.Range("E47").Value = "  -0.54%  "

$ws.Range("E48").Value = "  -1.21%  "

$ws.Range("E49").Value = "  +5.08%  "

$ws.Range("E50").Value = "  -0.54%  "

$c = $ws.Range("D51")
$c.NumberFormat = "@"
$c.Value = "0.998"
$c.Style = "Normal"
$ws.Range("E51").Value = "  -3.88%  "
